# [PHOENIX-5880] Modified Miscellaneous receipt scenario
# Adds a new "paymentMethod" worksheet (test data for payment-method driven
# miscellaneous receipts) after the existing "approvalDetails" sheet, makes
# it the active/selected sheet, and nudges the column widths on the other
# sheets to reflect the re-layout that happened when the workbook was
# re-saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "paymentMethod" sheet as the last tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "paymentMethod"

# Header row
$ws4.Range("A1").Value = "dataId"
$ws4.Range("B1").Value = "dd/chequeNum"
$ws4.Range("C1").Value = "bankName"

# Row 2 - cash
$ws4.Range("A2").Value = "cash"
$ws4.Range("B2").Value = "null"
$ws4.Range("C2").Value = "null"

# Row 3 - cheque
$ws4.Range("A3").Value = "cheque"
$ws4.Range("B3").Value = 123456
$ws4.Range("C3").Value = 102

# Row 4 - dd
$ws4.Range("A4").Value = "dd"
$ws4.Range("B4").Value = 123456
$ws4.Range("C4").Value = 102

# Row 5 - credit/debit card
$ws4.Range("A5").Value = "credit/debit card"

# Row 6 - direct bank
$ws4.Range("A6").Value = "direct bank"

# Row heights to match the re-flowed layout
$ws4.Rows.Item(1).RowHeight = 12.8
$ws4.Rows.Item(2).RowHeight = 20.2
$ws4.Rows.Item(3).RowHeight = 17.2
$ws4.Rows.Item(4).RowHeight = 19.45
$ws4.Rows.Item(5).RowHeight = 17.2
$ws4.Rows.Item(6).RowHeight = 20.2

# Column widths for the new sheet
$ws4.Columns.Item(1).ColumnWidth = 15.0407407407408
$ws4.Columns.Item(2).ColumnWidth = 12.8851851851852
$ws4.Columns.Item(3).ColumnWidth = 11.9074074074074

# Make the new sheet the active tab / selection
$null = $ws4.Activate()
$null = $ws4.Range("C4").Select()

# ---------------------------------------------------------------------
# 2. challanHeaderDetails: scroll view back to A1 (was F1).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("challanHeaderDetails")
$null = $ws2.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Column width adjustments on the other sheets (re-layout side effect
#    of the edit).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("chequeDetails")
$ws1.Columns.Item(1).ColumnWidth = 22.1962962962963
$ws1.Columns.Item(2).ColumnWidth = 16.1185185185186
$ws1.Columns.Item(3).ColumnWidth = 8.37777777777778
$ws1.Columns.Item(4).ColumnWidth = 25.2333333333334

$ws2.Columns.Item(1).ColumnWidth = 8.96666666666667
$ws2.Columns.Item(2).ColumnWidth = 18.862962962963
$ws2.Columns.Item(3).ColumnWidth = 19.3518518518519
$ws2.Columns.Item(4).ColumnWidth = 18.6666666666667
$ws2.Columns.Item(5).ColumnWidth = 18.6666666666667
$ws2.Columns.Item(6).ColumnWidth = 19.3518518518519
$ws2.Columns.Item(7).ColumnWidth = 19.3518518518519
$ws2.Columns.Item(8).ColumnWidth = 8.96666666666667

$ws3 = $wb.Worksheets.Item("approvalDetails")
$ws3.Columns.Item(1).ColumnWidth = 18.962962962963
$ws3.Columns.Item(2).ColumnWidth = 29.1518518518519
$ws3.Columns.Item(3).ColumnWidth = 29.4481481481482
$ws3.Columns.Item(4).ColumnWidth = 46.4962962962963

# ---------------------------------------------------------------------
# 4. Re-select the new sheet/cell last so it is the active tab on save
#    (activeTab moves from approvalDetails (2) to paymentMethod (3)).
# ---------------------------------------------------------------------
$null = $ws4.Activate()
$null = $ws4.Range("C4").Select()
